$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped by
# one day (45202 -> 45203) for every data row (rows 2 through 459).
$ws.Range("C2:C459").Value = 45203
